# Update "想去人数" (interested-count) figures in both the "展览" sheet
# and the duplicated "全部类型" sheet, as published at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1575
    3  = 49
    4  = 1028
    5  = 26
    7  = 2654
    9  = 1676
    10 = 183
    11 = 69
    14 = 10
    15 = 56
    16 = 72
    18 = 10
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates[$row]
}
$ws1.Cells.Item(12, 6).Value = 560

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates[$row]
}
$ws4.Cells.Item(12, 6).Value = 561
